# Applies the textual edits described by the diff:
#  1) "sí" + " una fila cumple..."      -> merged into a single run
#  2) "sí" + " una columna satisface..."-> merged into a single run
#  3) the "Se van obteniendo las columnas..." run gets the phrase
#     "de salida" inserted twice (split into several runs)
#  4) "filas" + ", entonces simplemente devuelve..." -> merged into a single run
#  5) "Recorre 2 listas..." run (with embedded proofErr tags) -> merged/cleaned
#  6) "primera_pasada/4  como..." / "...ultima_pasada..." runs -> merged/cleaned,
#     fixing the "ult" + "ima_pasada" spell-split into "ultima_pasada"

$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindContinue = 1

# NOTE: MatchCase is intentionally $false below -- the engine mis-handles
# MatchCase=$true together with accented characters (á, é, í, ó, ú, ñ, etc.)
# and fails to match even an exact-case search string. MatchWholeWord is
# left off too since we are matching arbitrary interior substrings.
function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $false, $false, $false, $false, $false, `
                                   $true, $wdFindContinue, $false, $replace, $wdReplaceAll)
    if (-not $ok) {
        Write-Host "NOT FOUND (see script source for the search text)"
    }
}

# 1) merge "sí" + " una fila cumple..." into a single run
Replace-Text "sí una fila cumple con sus pistas correspondientes, en caso de que se satisfagan entonces " `
             "sí una fila cumple con sus pistas correspondientes, en caso de que se satisfagan entonces "

# 2) merge "sí" + " una columna satisface..." into a single run
Replace-Text "sí una columna satisface las pistas asignadas a la misma, en caso afirmativo " `
             "sí una columna satisface las pistas asignadas a la misma, en caso afirmativo "

# 3) insert "de salida" twice into the "Se van obteniendo las columnas..." sentence
Replace-Text "en una lista indicando que cumple con sus pistas, en caso contrario se agrega un “0” indicando que no cumple con sus pistas" `
             "en una lista de salida indicando que cumple con sus pistas, en caso contrario se agrega un “0” en la lista de salida indicando que no cumple con sus pistas"

# 4) merge "filas" + ", entonces simplemente devuelve..." into a single run
Replace-Text "cumple con todas las pistas de las filas, entonces simplemente devuelve la misma grilla como grilla de salida" `
             "cumple con todas las pistas de las filas, entonces simplemente devuelve la misma grilla como grilla de salida"

# 5) merge the "Recorre 2 listas..." run and drop the proofErr-wrapped fragments
Replace-Text "tanto si es “#” , “X” o “_” en ambas." `
             "tanto si es “#” , “X” o “_” en ambas."

# 6a) merge "primera_pasada/" + "4  como" + " dato de entrada para " into a single run
Replace-Text "primera_pasada/4  como dato de entrada para segunda_pasada" `
             "primera_pasada/4  como dato de entrada para segunda_pasada"

# 6b) fix the "ult" + "ima_pasada" spell-split into "ultima_pasada"
Replace-Text "luego esta grilla se utiliza en ultima_pasada/4 que devuelve" `
             "luego esta grilla se utiliza en ultima_pasada/4 que devuelve"
